$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of row 301 down through the new rows (302:328)
# so column A keeps the date style (s="2") matching the existing data rows.
$ws.Range("A301:D301").Copy()
$ws.Range("A302:D328").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New daily data rows (dates 2021-06-29 .. 2021-07-25 -> serials 44376..44402)
$newRows = @(
  @(302, 44376, 0, 0, 0),
  @(303, 44377, 0, 0, 0),
  @(304, 44378, 0, 0, 0),
  @(305, 44379, 0, 0, 0),
  @(306, 44380, 0, 0, 0),
  @(307, 44381, 0, 0, 0),
  @(308, 44382, 0, 0, 0),
  @(309, 44383, 2, 2, 28.44950213371266),
  @(310, 44384, 0, 2, 28.44950213371266),
  @(311, 44385, 0, 2, 28.44950213371266),
  @(312, 44386, 0, 2, 28.44950213371266),
  @(313, 44387, 0, 2, 28.44950213371266),
  @(314, 44388, 0, 2, 28.44950213371266),
  @(315, 44389, 0, 2, 28.44950213371266),
  @(316, 44390, 0, 0, 0),
  @(317, 44391, 0, 0, 0),
  @(318, 44392, 0, 0, 0),
  @(319, 44393, 0, 0, 0),
  @(320, 44394, 0, 0, 0),
  @(321, 44395, 0, 0, 0),
  @(322, 44396, 0, 0, 0),
  @(323, 44397, 0, 0, 0),
  @(324, 44398, 0, 0, 0),
  @(325, 44399, 0, 0, 0),
  @(326, 44400, 0, 0, 0),
  @(327, 44401, 0, 0, 0),
  @(328, 44402, 1, 1, 14.22475106685633)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
